$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns B (lamda_1), C (lamda_2), D (dic_nbre_clients_poisson_2_keys),
# E (dic_nbre_clients_prob_poisson_2_values) for rows 2..55, matching the updated
# "auto scs and time in ms and auto capacity" figures.
$data = @(
    @(33.94444444444444, 1.95, 0, 0.162),
    @(33.94444444444444, 1.95, 3, 0.002),
    @(33.94444444444444, 1.95, 4, 0.006),
    @(33.94444444444444, 1.95, 5, 0.021),
    @(33.94444444444444, 1.95, 6, 0.028),
    @(33.94444444444444, 1.95, 7, 0.052),
    @(33.94444444444444, 1.95, 8, 0.044),
    @(33.94444444444444, 1.95, 9, 0.038),
    @(33.94444444444444, 1.95, 10, 0.033),
    @(33.94444444444444, 1.95, 11, 0.023),
    @(33.94444444444444, 1.95, 12, 0.033),
    @(33.94444444444444, 1.95, 13, 0.027),
    @(33.94444444444444, 1.95, 14, 0.023),
    @(33.94444444444444, 1.95, 15, 0.032),
    @(33.94444444444444, 1.95, 16, 0.03),
    @(33.94444444444444, 1.95, 17, 0.034),
    @(33.94444444444444, 1.95, 18, 0.03),
    @(33.94444444444444, 1.95, 19, 0.03),
    @(33.94444444444444, 1.95, 20, 0.027),
    @(33.94444444444444, 1.95, 21, 0.035),
    @(33.94444444444444, 1.95, 22, 0.027),
    @(33.94444444444444, 1.95, 23, 0.016),
    @(33.94444444444444, 1.95, 24, 0.028),
    @(33.94444444444444, 1.95, 25, 0.019),
    @(33.94444444444444, 1.95, 26, 0.016),
    @(33.94444444444444, 1.95, 27, 0.017),
    @(33.94444444444444, 1.95, 28, 0.021),
    @(33.94444444444444, 1.95, 29, 0.021),
    @(33.94444444444444, 1.95, 30, 0.014),
    @(33.94444444444444, 1.95, 31, 0.018),
    @(33.94444444444444, 1.95, 32, 0.011),
    @(33.94444444444444, 1.95, 33, 0.006),
    @(33.94444444444444, 1.95, 34, 0.003),
    @(33.94444444444444, 1.95, 35, 0.004),
    @(33.94444444444444, 1.95, 36, 0.014),
    @(33.94444444444444, 1.95, 37, 0.011),
    @(33.94444444444444, 1.95, 38, 0.006),
    @(33.94444444444444, 1.95, 39, 0.004),
    @(33.94444444444444, 1.95, 40, 0.004),
    @(33.94444444444444, 1.95, 41, 0.008),
    @(33.94444444444444, 1.95, 42, 0.002),
    @(33.94444444444444, 1.95, 43, 0.002),
    @(33.94444444444444, 1.95, 44, 0.003),
    @(33.94444444444444, 1.95, 45, 0.003),
    @(33.94444444444444, 1.95, 46, 0.001),
    @(33.94444444444444, 1.95, 47, 0.002),
    @(33.94444444444444, 1.95, 50, 0.001),
    @(33.94444444444444, 1.95, 54, 0.001),
    @(33.94444444444444, 1.95, 55, 0.001),
    @(33.94444444444444, 1.95, 56, 0.001),
    @(33.94444444444444, 1.95, 57, 0.001),
    @(33.94444444444444, 1.95, 60, 0.001),
    @(33.94444444444444, 1.95, 61, 0.001),
    @(33.94444444444444, 1.95, 71, 0.001)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $row = $row + 1
}
